$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 1813.625
$ws.Range("J2").Value = 5599.5
$ws.Range("L2").Value = 5599.5
$ws.Range("N2").Value = -5825.5
$ws.Range("H33").Value = 186.85715
$ws.Range("I33").Value = 190.54546
$ws.Range("J33").Value = 173.33333
$ws.Range("K33").Value = 190.54546
$ws.Range("L33").Value = 173.33333
$ws.Range("M33").Value = 38.45454000000001
$ws.Range("N33").Value = -631.3333299999999
$ws.Range("H53").Value = 336.41666
$ws.Range("I53").Value = 315.1111
$ws.Range("J53").Value = 400.33334
$ws.Range("K53").Value = 315.1111
$ws.Range("L53").Value = 400.33334
$ws.Range("M53").Value = 321.8889
$ws.Range("N53").Value = -1674.33334
$ws.Range("H64").Value = 4999.6665
$ws.Range("I64").Value = 4999.6665
$ws.Range("K64").Value = 4999.6665
$ws.Range("M64").Value = -4751.6665
$ws.Range("H67").Value = 4999.6665
$ws.Range("I67").Value = 4999.6665
$ws.Range("K67").Value = 4999.6665
$ws.Range("M67").Value = -4141.6665
$ws.Range("H80").Value = 728.25
$ws.Range("I80").Value = 1150
$ws.Range("J80").Value = 475.2
$ws.Range("K80").Value = 3450
$ws.Range("L80").Value = 1425.6
$ws.Range("M80").Value = -2452
$ws.Range("N80").Value = -3421.6
$ws.Range("H83").Value = 728.25
$ws.Range("I83").Value = 1150
$ws.Range("J83").Value = 475.2
$ws.Range("K83").Value = 10350
$ws.Range("L83").Value = 4276.8
$ws.Range("M83").Value = -5358
$ws.Range("N83").Value = -14260.8

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2450
$ws.Range("I61").Value = 2345.5
$ws.Range("K61").Value = 2345.5
$ws.Range("M61").Value = -2133.5
$ws.Range("H136").Value = 2450
$ws.Range("I136").Value = 2345.5
$ws.Range("K136").Value = 7036.5
$ws.Range("M136").Value = -4486.5
$ws.Range("H137").Value = 60000

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 6371.9165
$ws.Range("I20").Value = 6215.75
$ws.Range("J20").Value = 6684.25
$ws.Range("K20").Value = 6215.75
$ws.Range("L20").Value = 6684.25
$ws.Range("M20").Value = -5968.75
$ws.Range("N20").Value = -7178.25
$ws.Range("H86").Value = 3373
$ws.Range("I86").Value = 3247.5
$ws.Range("J86").Value = 3498.5
$ws.Range("K86").Value = 3247.5
$ws.Range("L86").Value = 3498.5
$ws.Range("M86").Value = -2124.5
$ws.Range("N86").Value = -5744.5
$ws.Range("H89").Value = 3373
$ws.Range("I89").Value = 3247.5
$ws.Range("J89").Value = 3498.5
$ws.Range("K89").Value = 16237.5
$ws.Range("L89").Value = 17492.5
$ws.Range("M89").Value = -10621.5
$ws.Range("N89").Value = -28724.5
$ws.Range("H130").Value = 44899.5
$ws.Range("J130").Value = 44899.5
$ws.Range("L130").Value = 44899.5
$ws.Range("N130").Value = -54939.5

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 49
$ws.Range("I7").Value = 54.285713
$ws.Range("J7").Value = 39.75
$ws.Range("K7").Value = 54.285713
$ws.Range("L7").Value = 39.75
$ws.Range("M7").Value = 58.714287
$ws.Range("N7").Value = -265.75
$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("M22").ClearContents()
$ws.Range("H86").Value = 14933
$ws.Range("H89").Value = 14933

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H61").Value = 357.6
$ws.Range("I61").Value = 415.8
$ws.Range("J61").Value = 328.5
$ws.Range("K61").Value = 1247.4
$ws.Range("L61").Value = 985.5
$ws.Range("M61").Value = -1032.4
$ws.Range("N61").Value = -1415.5
$ws.Range("H75").Value = 2866
$ws.Range("J75").Value = 2902.6
$ws.Range("L75").Value = 8707.799999999999
$ws.Range("N75").Value = -10703.8
$ws.Range("H78").Value = 2866
$ws.Range("J78").Value = 2902.6
$ws.Range("L78").Value = 26123.4
$ws.Range("N78").Value = -36107.39999999999
$ws.Range("H103").Value = 184.88889
$ws.Range("I103").Value = 196.57143
$ws.Range("J103").Value = 144
$ws.Range("K103").Value = 589.71429
$ws.Range("L103").Value = 432
$ws.Range("M103").Value = 289.28571
$ws.Range("N103").Value = -2190
$ws.Range("H114").Value = 339.83334
$ws.Range("J114").Value = 400.25
$ws.Range("L114").Value = 1200.75
$ws.Range("N114").Value = -7708.75
$ws.Range("H117").Value = 756
$ws.Range("I117").Value = 756
$ws.Range("K117").Value = 2268
$ws.Range("M117").Value = 1174
$ws.Range("H129").Value = 622.2857
$ws.Range("I129").Value = 585
$ws.Range("J129").Value = 715.5
$ws.Range("K129").Value = 1755
$ws.Range("L129").Value = 2146.5
$ws.Range("M129").Value = 3245
$ws.Range("N129").Value = -12146.5

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 2414.8823
$ws.Range("I122").Value = 1526.8462
$ws.Range("J122").Value = 5301
$ws.Range("K122").Value = 4580.5386
$ws.Range("L122").Value = 15903
$ws.Range("M122").Value = -2130.5386
$ws.Range("N122").Value = -20803

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 0
$ws.Range("I7").Value = 0
$ws.Range("K7").Value = 0
$ws.Range("M7").ClearContents()
$ws.Range("H22").Value = 1508.7273
$ws.Range("I22").Value = 1349.25
$ws.Range("J22").Value = 1934
$ws.Range("K22").Value = 1349.25
$ws.Range("L22").Value = 1934
$ws.Range("M22").Value = -1054.25
$ws.Range("N22").Value = -2524
$ws.Range("H27").Value = 1508.7273
$ws.Range("I27").Value = 1349.25
$ws.Range("J27").Value = 1934
$ws.Range("K27").Value = 1349.25
$ws.Range("L27").Value = 1934
$ws.Range("M27").Value = -1242.25
$ws.Range("N27").Value = -2148
$ws.Range("H55").Value = 311.55554
$ws.Range("I55").Value = 265.6
$ws.Range("J55").Value = 369
$ws.Range("K55").Value = 265.6
$ws.Range("L55").Value = 369
$ws.Range("M55").Value = -92.60000000000002
$ws.Range("N55").Value = -715
$ws.Range("H56").Value = 21987.75
$ws.Range("I56").Value = 21987.75
$ws.Range("K56").Value = 21987.75
$ws.Range("M56").Value = -21296.75
$ws.Range("H76").Value = 35000
$ws.Range("J76").Value = 35000
$ws.Range("L76").Value = 35000
$ws.Range("N76").Value = -35676
$ws.Range("H79").Value = 35000
$ws.Range("J79").Value = 35000
$ws.Range("L79").Value = 35000
$ws.Range("N79").Value = -37340
$ws.Range("H106").Value = 8000
$ws.Range("J106").Value = 8000
$ws.Range("L106").Value = 8000
$ws.Range("N106").Value = -10524
$ws.Range("H126").Value = 0
$ws.Range("I126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("M126").ClearContents()

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H28").Value = 8000
$ws.Range("I28").Value = 0
$ws.Range("J28").Value = 8000
$ws.Range("K28").Value = 0
$ws.Range("L28").Value = 8000
$ws.Range("M28").ClearContents()
$ws.Range("N28").Value = -8696
$ws.Range("H69").Value = 19290.334
$ws.Range("J69").Value = 19290.334
$ws.Range("L69").Value = 19290.334
$ws.Range("N69").Value = -20788.334
$ws.Range("H72").Value = 19290.334
$ws.Range("J72").Value = 19290.334
$ws.Range("L72").Value = 57871.00199999999
$ws.Range("N72").Value = -65359.00199999999
